# "update baseline and EE scenariio"
#
# Source data behind the "EE Trajectory" sheet's 2030 column (S41:S45) was
# re-pointed from the 2025 figures (column AB) to the 2030 figures
# (column AC) on the little lookup block at Y3:AE8 on the same sheet.
# That ripples through the rest of the interpolated trajectory row
# (C:S uses the 2013-2030 step, T:AM uses the 2030-2050 step) and the
# row-46 totals.
#
# The user also ended their session on the "EE Trajectory" tab instead of
# "Baseline Trajectory", with a different cell selected.

$wb = $excel.ActiveWorkbook

$wsBaseline = $wb.Worksheets.Item("Baseline Trajectory")
$wsEE       = $wb.Worksheets.Item("EE Trajectory")

# --- content edit: re-point the 2030 baseline-year lookups from the 2025
# column (AB) to the 2030 column (AC) ---
$wsEE.Range("S41").Formula    = "=AC4"
$wsEE.Range("S42:S45").Formula = "=AC5"

# --- view/selection: EE Trajectory becomes the active sheet/tab, with a
# new selected cell; Baseline Trajectory keeps its prior selection ---
$wsEE.Activate() | Out-Null
$wsEE.Range("AO42").Select() | Out-Null
